# Generate Report for handoff
# Adds a new localization entry (9e3ee29f-...) that is now "Ready for handoff"
# and moves the previously-"Ready for handoff" files (8b3f1882-..., d8590523-...)
# into "In Translation" status. Also records a new handoff for e398423b-...
# and shifts the ".localization-config" bookkeeping row down.

$wb = $excel.ActiveWorkbook

function Add-Or-UpdateHyperlink($ws, $cellAddr, $url, $displayText) {
    $ws.Hyperlinks.Add($ws.Range($cellAddr), $url, [Type]::Missing, [Type]::Missing, $displayText) | Out-Null
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Existing rows 2 & 3 now show "In Translation" instead of "Ready for handoff"
$wsOverview.Range("B2").Value = "In Translation"
$wsOverview.Range("C2").Value = "In Translation"
$wsOverview.Range("B3").Value = "In Translation"
$wsOverview.Range("C3").Value = "In Translation"

# Row 4 now refers to the newly added file 9e3ee29f-... (status "Ready for handoff")
$wsOverview.Range("A4").Value = "9e3ee29f-568b-4524-9606-539e51e14127.md"
$wsOverview.Range("B4").Value = "Ready for handoff"
$wsOverview.Range("C4").Value = "Ready for handoff"
Add-Or-UpdateHyperlink $wsOverview "A4" "https://github.com/OpenLocalizationTest/oltest/blob/b6a04f1d45a8721021730a766a2d7b9a8ae54e15/e2e/9e3ee29f-568b-4524-9606-539e51e14127.md" "9e3ee29f-568b-4524-9606-539e51e14127.md"

# New row 5: e398423b-... (status "Ready for handoff")
$wsOverview.Range("A5").Value = "e398423b-5319-4e64-9d0a-3a03345725a5.md"
$wsOverview.Range("B5").Value = "Ready for handoff"
$wsOverview.Range("C5").Value = "Ready for handoff"
Add-Or-UpdateHyperlink $wsOverview "A5" "https://github.com/OpenLocalizationTest/oltest/blob/b6a04f1d45a8721021730a766a2d7b9a8ae54e15/e2e/e398423b-5319-4e64-9d0a-3a03345725a5.md" "e398423b-5319-4e64-9d0a-3a03345725a5.md"

# New row 6: the ".localization-config" bookkeeping row, now pushed down to row 6
$wsOverview.Range("A6").Value = ".localization-config"
$wsOverview.Range("B6").Value = "Not to be localized"
$wsOverview.Range("C6").Value = "Not to be localized"
Add-Or-UpdateHyperlink $wsOverview "A6" "https://github.com/OpenLocalizationTest/oltest/blob/b6a04f1d45a8721021730a766a2d7b9a8ae54e15/.localization-config" ".localization-config"

# Copy formatting (font / number format) from existing rows onto the new rows so the
# new cells look consistent with the rest of the table (hyperlink style for column A).
$wsOverview.Range("A3").Copy() | Out-Null
$wsOverview.Range("A4:A6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$wsOverview.Range("B4").Copy() | Out-Null
$wsOverview.Range("B5:C6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("B2").Value = "In Translation"
$wsZh.Range("B3").Value = "In Translation"

# Row 4 now refers to 9e3ee29f-...
$wsZh.Range("A4").Value = "9e3ee29f-568b-4524-9606-539e51e14127.md"
$wsZh.Range("B4").Value = "Ready for handoff"
$wsZh.Range("C4").Value = "9e3ee29f-568b-4524-9606-539e51e14127.9a73fe992e06dc0220eb32707729b5e5e3da8da7.zh-cn.xlf"
$wsZh.Range("D4").Value = "2016-01-19 05:06:55"
$wsZh.Range("G4").Value = "0001-01-01 00:00:00"
$wsZh.Range("H4").Value = "Include"
Add-Or-UpdateHyperlink $wsZh "A4" "https://github.com/OpenLocalizationTest/oltest/blob/b6a04f1d45a8721021730a766a2d7b9a8ae54e15/e2e/9e3ee29f-568b-4524-9606-539e51e14127.md" "9e3ee29f-568b-4524-9606-539e51e14127.md"
Add-Or-UpdateHyperlink $wsZh "C4" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7485f37e86a6767d70a5def5900399e1124712dd/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/9e3ee29f-568b-4524-9606-539e51e14127.9a73fe992e06dc0220eb32707729b5e5e3da8da7.zh-cn.xlf" "9e3ee29f-568b-4524-9606-539e51e14127.9a73fe992e06dc0220eb32707729b5e5e3da8da7.zh-cn.xlf"

# New row 5: e398423b-...
$wsZh.Range("A5").Value = "e398423b-5319-4e64-9d0a-3a03345725a5.md"
$wsZh.Range("B5").Value = "Ready for handoff"
$wsZh.Range("C5").Value = "e398423b-5319-4e64-9d0a-3a03345725a5.92e40c4cbd82d1cdc2115d9aa5242ed4622f271a.zh-cn.xlf"
$wsZh.Range("D5").Value = "2016-01-19 05:06:55"
$wsZh.Range("G5").Value = "0001-01-01 00:00:00"
$wsZh.Range("H5").Value = "Include"
Add-Or-UpdateHyperlink $wsZh "A5" "https://github.com/OpenLocalizationTest/oltest/blob/b6a04f1d45a8721021730a766a2d7b9a8ae54e15/e2e/e398423b-5319-4e64-9d0a-3a03345725a5.md" "e398423b-5319-4e64-9d0a-3a03345725a5.md"
Add-Or-UpdateHyperlink $wsZh "C5" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7485f37e86a6767d70a5def5900399e1124712dd/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/e398423b-5319-4e64-9d0a-3a03345725a5.92e40c4cbd82d1cdc2115d9aa5242ed4622f271a.zh-cn.xlf" "e398423b-5319-4e64-9d0a-3a03345725a5.92e40c4cbd82d1cdc2115d9aa5242ed4622f271a.zh-cn.xlf"

# New row 6: ".localization-config" bookkeeping row
$wsZh.Range("A6").Value = ".localization-config"
$wsZh.Range("B6").Value = "Not to be localized"
$wsZh.Range("D6").Value = "0001-01-01 00:00:00"
$wsZh.Range("G6").Value = "0001-01-01 00:00:00"
$wsZh.Range("H6").Value = "Ignored"
Add-Or-UpdateHyperlink $wsZh "A6" "https://github.com/OpenLocalizationTest/oltest/blob/b6a04f1d45a8721021730a766a2d7b9a8ae54e15/.localization-config" ".localization-config"

# Formatting for new rows (column styles copied from row 3, which already has the
# full set of styles: hyperlink cells in A/C, plain cells in B, date-styled D/G/H)
$wsZh.Range("A3:I3").Copy() | Out-Null
$wsZh.Range("A4:I4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$wsZh.Range("A3:I3").Copy() | Out-Null
$wsZh.Range("A5:I5").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$wsZh.Range("A3:I3").Copy() | Out-Null
$wsZh.Range("A6:I6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B2").Value = "In Translation"
$wsDe.Range("B3").Value = "In Translation"

# Row 4 now refers to 9e3ee29f-...
$wsDe.Range("A4").Value = "9e3ee29f-568b-4524-9606-539e51e14127.md"
$wsDe.Range("B4").Value = "Ready for handoff"
$wsDe.Range("C4").Value = "9e3ee29f-568b-4524-9606-539e51e14127.9a73fe992e06dc0220eb32707729b5e5e3da8da7.de-de.xlf"
$wsDe.Range("D4").Value = "2016-01-19 05:07:04"
$wsDe.Range("G4").Value = "0001-01-01 00:00:00"
$wsDe.Range("H4").Value = "Include"
Add-Or-UpdateHyperlink $wsDe "A4" "https://github.com/OpenLocalizationTest/oltest/blob/b6a04f1d45a8721021730a766a2d7b9a8ae54e15/e2e/9e3ee29f-568b-4524-9606-539e51e14127.md" "9e3ee29f-568b-4524-9606-539e51e14127.md"
Add-Or-UpdateHyperlink $wsDe "C4" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1824d9773e06f0be8a256a095409b58e3df7149e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/9e3ee29f-568b-4524-9606-539e51e14127.9a73fe992e06dc0220eb32707729b5e5e3da8da7.de-de.xlf" "9e3ee29f-568b-4524-9606-539e51e14127.9a73fe992e06dc0220eb32707729b5e5e3da8da7.de-de.xlf"

# New row 5: e398423b-...
$wsDe.Range("A5").Value = "e398423b-5319-4e64-9d0a-3a03345725a5.md"
$wsDe.Range("B5").Value = "Ready for handoff"
$wsDe.Range("C5").Value = "e398423b-5319-4e64-9d0a-3a03345725a5.92e40c4cbd82d1cdc2115d9aa5242ed4622f271a.de-de.xlf"
$wsDe.Range("D5").Value = "2016-01-19 05:07:04"
$wsDe.Range("G5").Value = "0001-01-01 00:00:00"
$wsDe.Range("H5").Value = "Include"
Add-Or-UpdateHyperlink $wsDe "A5" "https://github.com/OpenLocalizationTest/oltest/blob/b6a04f1d45a8721021730a766a2d7b9a8ae54e15/e2e/e398423b-5319-4e64-9d0a-3a03345725a5.md" "e398423b-5319-4e64-9d0a-3a03345725a5.md"
Add-Or-UpdateHyperlink $wsDe "C5" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1824d9773e06f0be8a256a095409b58e3df7149e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/e398423b-5319-4e64-9d0a-3a03345725a5.92e40c4cbd82d1cdc2115d9aa5242ed4622f271a.de-de.xlf" "e398423b-5319-4e64-9d0a-3a03345725a5.92e40c4cbd82d1cdc2115d9aa5242ed4622f271a.de-de.xlf"

# New row 6: ".localization-config" bookkeeping row
$wsDe.Range("A6").Value = ".localization-config"
$wsDe.Range("B6").Value = "Not to be localized"
$wsDe.Range("D6").Value = "0001-01-01 00:00:00"
$wsDe.Range("G6").Value = "0001-01-01 00:00:00"
$wsDe.Range("H6").Value = "Ignored"
Add-Or-UpdateHyperlink $wsDe "A6" "https://github.com/OpenLocalizationTest/oltest/blob/b6a04f1d45a8721021730a766a2d7b9a8ae54e15/.localization-config" ".localization-config"

# Formatting for new rows
$wsDe.Range("A3:I3").Copy() | Out-Null
$wsDe.Range("A4:I4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$wsDe.Range("A3:I3").Copy() | Out-Null
$wsDe.Range("A5:I5").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$wsDe.Range("A3:I3").Copy() | Out-Null
$wsDe.Range("A6:I6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

$wb.Save()
